$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Finalize "Vista de detalle de compras" (purchases detail view) ---

# Mark the "Cliente" (row 7) and "caja" (row 17) tasks as reviewed by
# highlighting their rows in yellow, matching the rest of the finished rows.
$ws.Range("D7:H7").Interior.Color = 65535
$ws.Range("D17:H17").Interior.Color = 65535

# "Orden_Compra" (row 11) and "det_compra" (row 12) purchase tasks are now
# complete: mark them compliant and at 100%.
$ws.Range("G11").Value = "si"
$ws.Range("H11").Value = 100
$ws.Range("G12").Value = "si"
$ws.Range("H12").Value = 100

# Reveal the detail columns (Fecha Entrega / Cumple / % Cumplimiento) that
# were previously hidden, now that the purchases detail is finished.
$ws.Columns("F").Hidden = $false
$ws.Columns("G").Hidden = $false
$ws.Columns("G").ColumnWidth = 7.85546875
$ws.Columns("H").Hidden = $false

# Reset the view: scroll back to the top and move the selection to H20.
$ws.Range("H20").Select()
